$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -ne $val) {
        $s = [string]$val
        if ($s.EndsWith("16")) {
            $cell.Value2 = $s.Substring(0, $s.Length - 2)
        }
    }
}
